# Agregado de etiquetas menudencias: lengua, molleja, quijada, tendon - jabat
# Completa las columnas M (codigo) y N (estado) para las filas 35-38 de Hoja1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fila 35 - Lengua
$ws.Range("M35").Value = 161
$ws.Range("N35").Value = "OK"

# Fila 36 - Quijada
$ws.Range("M36").Value = 163
$ws.Range("N36").Value = "OK"

# Fila 37 - Tendon de brazuelo
$ws.Range("M37").Value = 177
$ws.Range("N37").Value = "OK"

# Fila 38 - Molleja
$ws.Range("M38").Value = 162
$ws.Range("N38").Value = "OK"

# Actualiza la posicion de desplazamiento/seleccion de la vista
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("N39").Select()
